$d = $word.ActiveDocument

# Locate the paragraph "The name of the app should be displayed correctly with
# proper font and color." - the new bullet needs to be inserted right after it
# (and before the "Compatibility Test Cases" heading).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "The name of the app should be displayed correctly with proper font and color.*") {
        $target = $p
        break
    }
}

# Split off a new, empty paragraph right after the target one. Word
# automatically carries over the ListParagraph style / numbering (numId 2)
# from the paragraph it was split from, so the new paragraph already has the
# correct <w:pPr> (pStyle=ListParagraph, numPr ilvl=0/numId=2).
$target.Range.InsertParagraphAfter() | Out-Null
$newRange = $target.Next().Range

# Fill the new paragraph with the runs from the diff, including the
# proofErr spell-check markers bracketing the two "Niki" occurrences and the
# curly opening quotation marks used in the source text.
$quote = [char]0x201C

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
</w:pPr>
<w:r><w:t xml:space="preserve">$quote </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Niki</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> is typing $quote text appears when </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Niki</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> is typing.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$newRange.InsertXML($xml) | Out-Null

$d.Save()
